$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.507.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.26%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.67%  "

# Row 6
$ws.Range("E6").Value = "  +4.26%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.73%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.342"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.02%  "

# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.82%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0710"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.91%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.71%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.180.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.44%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.956.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.14%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.699"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.17%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.87%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.512.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "

# Row 19
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.60%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0834"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.51%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "

# Row 22
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.15%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.99%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.81%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.89%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.95%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.962"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +21.88%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0572"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.53%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.90%  "

# Row 35
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.52%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.01%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "92.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0641"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.48%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0205"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.63%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.348.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.37%  "

# Row 45
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "49.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +41.00%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.06%  "

# Row 47
$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.37%  "

# Row 48
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "

# Row 50
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.75%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.091.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.48%  "
